$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column before column N (shifts old N,O,P -> O,P,Q).
# Excel copies formatting from the column to the left (M) for the new column.
$ws.Columns("N:N").Insert()

# The newly inserted column picks up the width of the column to its left (M).
$ws.Columns("N:N").ColumnWidth = $ws.Columns("M:M").ColumnWidth

# Switch the active sheet to "Repayment Schedule" and select cell P6,
# which also clears the previous active selection on "Transactions".
[void]$ws.Activate()
[void]$ws.Range("P6").Select()
